# Update the INVOICE_COUNTY (column T) empty placeholder values "--" to "MERKEZ"
# and bump the DATE_ADDED (column AA) values from 08/01/2024 to 09/01/2024
# across all data rows of the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).Worksheet.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $countyCell = $ws.Cells.Item($r, 20)
    if ($countyCell.Text -eq "--") {
        $countyCell.Value = "MERKEZ"
    }

    $dateCell = $ws.Cells.Item($r, 27)
    if ($dateCell.Text -eq "08/01/2024") {
        # Leading apostrophe forces Excel to keep this as literal text instead
        # of auto-converting the date-like string into a date serial number.
        $dateCell.Value = "'09/01/2024"
    }
}
